# Classwork as of 9/23
# Adds 1st/2nd/3rd quartile + IQR rows to the "Data" sheet, and tidies up
# the number-format/style noise that had accumulated on column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New statistics: 1st Quartile (E14), 2nd Quartile / median (E15),
#    3rd Quartile (E16) and Interquartile Range (E17).
# ---------------------------------------------------------------------
$ws.Range("E14").Formula = "=QUARTILE.EXC(B2:B13,1)"
$ws.Range("E15").Formula = "=PERCENTILE.INC(B2:B13,0.5)"
$ws.Range("E16").Formula = "=QUARTILE.EXC(B2:B13,3)"
$ws.Range("E17").Formula = "=E16-E14"

# ---------------------------------------------------------------------
# 2. Tidy up formatting on column E: cells that only ever displayed the
#    General number format had redundant "apply number format" styling;
#    nudging the font back onto itself collapses those cells onto the
#    shared (non currency) style instead of their own throw-away one.
# ---------------------------------------------------------------------
$plainCells = @("E1","E6","E7","E8","E9","E10","E11","E12","E13","E14","E15","E16","E17")
foreach ($addr in $plainCells) {
  $c = $ws.Range($addr)
  $c.Font.Name = "Times New Roman"
  $c.Font.Size = 12
  $c.Font.Bold = $false
}

# ---------------------------------------------------------------------
# 3. Move the active selection to I12 (matches the saved cursor position).
# ---------------------------------------------------------------------
$ws.Range("I12").Select()
